$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "The code is currently structured with 5 classes ..." paragraph rewrite
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The code is currently structured with 5 classes 2 being used as control classes (explained later) and the other three classes being used to store methods for each type of user or the Main class wich exists to control the menu and start the program along side call the other classes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The code is currently structured with 4 classes 2 being used as control classes (explained later) and the other two classes being used to store methods for each one being used to store variables of a users account and the main which is for launching the program and creating a connection to the database",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "The control classes exsist ..." paragraph: split into two paragraphs
#    and rewrite the tail text of the original + add the new second
#    paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The control classes exsist to hold functions that activate based on the JavaFX buttons and help control the flow of the GUI each javafx file contains code that when a button is pressed that it calls the appropriate control type with both user types having seprate control classes when a function is unshared to make the code more secure and readable ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The control classes exsist to hold functions that activate based on the JavaFX buttons and help control the flow of the GUI each javafx file contains code that when a button is pressed that it calls the appropriate control type with both user ^pIt also exists to hold functions that modfy the user control class a variable of current user is stored to see who is currently logged in which is read from the database and the administrator control class exsist to hold functions for administrator actions such as creating an account and loading administrator GUI",
    2) | Out-Null

# Two extra blank paragraphs are inserted right after the new second
# paragraph (before the pre-existing run of blank paragraphs / page break).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -like "It also exists to hold functions*administrator GUI") {
        $target = $i
        break
    }
}
$pAdmin = $d.Paragraphs($target)
$pAdmin.Range.InsertParagraphAfter()
$d.Paragraphs($target + 1).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3. Two extra blank paragraphs after the "Start shift and end shift- ..."
#    paragraph.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -like "Start shift and end shift-*") {
        $target = $i
        break
    }
}
$pShift = $d.Paragraphs($target)
$pShift.Range.InsertParagraphAfter()
$d.Paragraphs($target + 1).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 4. Replace the whole "Adminstartor class" section with a "User Controller
#    class" Heading 1, dropping all the paragraphs that used to sit between
#    it and the following page break.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Adminstartor class") {
        $target = $i
        break
    }
}

$pReadAccount = $null
for ($i = $target + 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -like "Read account*") {
        $pReadAccount = $i
        break
    }
}

# Delete everything from right after the heading through the end of the
# "Read account ..." paragraph first (so paragraph indices at/below the
# heading do not shift before we use them), then fix up the heading itself.
$startPos = $d.Paragraphs($target + 1).Range.Start
$endPos = $d.Paragraphs($pReadAccount).Range.End
$rngKill = $d.Range($startPos, $endPos)
$rngKill.Delete()

$pHeading = $d.Paragraphs($target)
$pHeading.Range.Text = "User Controller class"
$pHeading.Style = "Heading 1"

Write-Output "done"
